{"js": "// Fix a handful of typos / wording tweaks in the \"8.3 Final Project Step 1\" doc.\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"Probem Statment: How can we decearse waste (shrink) in the produce departemnt?\",\n    replace: \"Probem Statment: How can we decearse waste (shrink) in the produce department?\",\n  },\n  {\n    find: \"We have learned how to create many different visualizations in this course to help us better understnad our data and uncover trends/relationships. I think a few of the important ones that we will use in this analysis will be the below;\",\n    replace: \"We have learned how to create many different visualizations in this course to help us better understand our data and uncover trends/relationships. I think a few of the important ones that we will use in this analysis will be the below;\",\n  },\n  {\n    find: \"Scatter plots - for visualizing relationships.\",\n    replace: \"Scatter plots - for visualizing relationships between all variables.\",\n  },\n  {\n    find: \"Bar charts (Pedro chart) - for identifying top areas of focus/impact.\",\n    replace: \"Bar charts (Pedro chart) - for identifying top areas (specific produce fruit/vegetable) of focus/impact.\",\n  },\n  {\n    find: \"Line charts - for identifying seasonality\",\n    replace: \"Line charts - for identifying seasonality in sales data.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix a handful of typos / wording tweaks in the \"8.3 Final Project Step 1\" doc.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Probem Statment: How can we decearse waste (shrink) in the produce departemnt?\"; Replace = \"Probem Statment: How can we decearse waste (shrink) in the produce department?\" },\n    @{ Find = \"We have learned how to create many different visualizations in this course to help us better understnad our data and uncover trends/relationships. I think a few of the important ones that we will use in this analysis will be the below;\"; Replace = \"We have learned how to create many different visualizations in this course to help us better understand our data and uncover trends/relationships. I think a few of the important ones that we will use in this analysis will be the below;\" },\n    @{ Find = \"Scatter plots - for visualizing relationships.\"; Replace = \"Scatter plots - for visualizing relationships between all variables.\" },\n    @{ Find = \"Bar charts (Pedro chart) - for identifying top areas of focus/impact.\"; Replace = \"Bar charts (Pedro chart) - for identifying top areas (specific produce fruit/vegetable) of focus/impact.\" },\n    @{ Find = \"Line charts - for identifying seasonality\"; Replace = \"Line charts - for identifying seasonality in sales data.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.Find\n    $found = $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 1)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($r.Find)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
